# Rajout Fichier html character
# - add a new "Classe" worksheet after "Monster"
# - fill it with a header row (same headers as Monster, except A1 = "Classe")
#   and a single column of monster classes (Orque, Gobelin, UrukHaï, RoiSorcier, Troll)
# - tweak a couple of column widths on the new sheet
# - update the selection on the "Monster" sheet and move the active tab to "Classe"

$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "Monster" sheet --------------------
$wsMonster = $wb.Worksheets.Item("Monster")
$wsMonster.Range("B1:J1").Select()

# --- Add the new "Classe" sheet after the last existing sheet -----------
$sheetCount = $wb.Worksheets.Count
$wsClasse = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$wsClasse.Name = "Classe"

# --- Header row (mirrors Monster's header row, first column becomes "Classe") ---
$wsClasse.Range("A1").Value = "Classe"
$wsClasse.Range("B1").Value = "Strength"
$wsClasse.Range("C1").Value = "Hability"
$wsClasse.Range("D1").Value = "Endurance"
$wsClasse.Range("E1").Value = "Intelligence"
$wsClasse.Range("F1").Value = "Mana"
$wsClasse.Range("G1").Value = "Gain Xp"
$wsClasse.Range("H1").Value = "Level equivalent"
$wsClasse.Range("I1").Value = "Fréquence"
$wsClasse.Range("J1").Value = "Special habilities"

# --- Monster classes, column A ---
$wsClasse.Range("A2").Value = "Orque"
$wsClasse.Range("A3").Value = "Gobelin"
$wsClasse.Range("A4").Value = "UrukHaï"
$wsClasse.Range("A5").Value = "RoiSorcier"
$wsClasse.Range("A6").Value = "Troll"

# --- Column widths for columns H and J ---
$wsClasse.Columns.Item(8).ColumnWidth = 18.75
$wsClasse.Columns.Item(10).ColumnWidth = 18.75

# --- Page margins to match the rest of the workbook ---
$ps = $wsClasse.PageSetup
$ps.LeftMargin = 0.75 * 72
$ps.RightMargin = 0.75 * 72
$ps.TopMargin = 1 * 72
$ps.BottomMargin = 1 * 72
$ps.HeaderMargin = 0.5 * 72
$ps.FooterMargin = 0.5 * 72

# --- Selection / active cell on the new sheet ---
$wsClasse.Range("A12").Select()
